# Updates the "Estado de Cuenta" worksheet:
#  - refreshes the header totals (VALOR MORA, Cant. Trabajadores, Cant. Periodos)
#  - re-sorts the existing worker's period rows into ascending order
#  - appends the first new worker (part 1 of the new estados de cuenta)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header totals -----------------------------------------------------
$ws.Range("E11").Value = 449826
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 7

# --- make room for the two extra detail rows (new employee's row +
#     the table's closing/total-styled row) right after the existing
#     worker's last period row, pushing the signature block down.
$ws.Rows("22:23").Insert()

# Capture the "closing / totals" border formatting that currently lives
# on row 21 (thicker bottom border) before it gets overwritten below, so
# the new last row (23) can reuse it.
$ws.Range("B21:J21").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# --- re-sort CARLOS RAFAEL PANTOJA ROJAS's periods ascending (2502..2507)
#     row 16 keeps the reduced value (34164) but now represents period 2502
#     instead of 2507; rows 17-20 are the flat 56940 periods; row 21 moves
#     from the "totals" border style to the regular row style and becomes
#     period 2507 / 56940.
$ws.Range("B16:J16").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 34164

$ws.Range("E17").Value = "2503"
$ws.Range("F17").Value = 56940

$ws.Range("E18").Value = "2504"
$ws.Range("F18").Value = 56940

$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 56940

$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 56940

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1002191980"
$ws.Range("D21").Value = "CARLOS RAFAEL PANTOJA ROJAS"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# --- new worker: LUIS DANIEL AYALA ORTIZ, single period 2508 -----------
$ws.Range("B16:J16").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1143413951"
$ws.Range("D22").Value = "LUIS DANIEL AYALA ORTIZ"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 74022
$ws.Range("G22").Value = 1850550

# --- closing row of the table (same "totals" border style as the old
#     last row): repeats CARLOS's identity columns with period 2508
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1002191980"
$ws.Range("D23").Value = "CARLOS RAFAEL PANTOJA ROJAS"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500
